$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Thành tích học tập" (achievement) column (F) for each student row.
$ws.Range("F2").Value = "Học sinh giỏi"
$ws.Range("F3").Value = "Thành tích đặc biệt"
$ws.Range("F4").Value = "Không"
$ws.Range("F5").Value = "Học sinh tiên tiến"
$ws.Range("F6").Value = "Học sinh giỏi"
